# Update column G ("K") values for rows 2-15 on the active worksheet.
# These values were regenerated (K instead of Strike#) per the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6
    3  = 7
    4  = 3
    5  = 4
    6  = 6
    7  = 5
    8  = 8
    9  = 4
    10 = 6
    11 = 4
    12 = 4
    13 = 6
    14 = 4
    15 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
